$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "66÷5=13, 1"
$t.Cell(1, 2).Range.Text = "85÷8=10, 5"
$t.Cell(1, 3).Range.Text = "76÷6=12, 4"
$t.Cell(1, 4).Range.Text = "15÷7=2, 1"
$t.Cell(1, 5).Range.Text = "67÷4=16, 3"

$t.Cell(5, 1).Range.Text = "82÷6=13, 4"
$t.Cell(5, 2).Range.Text = "14÷5=2, 4"
$t.Cell(5, 3).Range.Text = "96÷3=32, 0"
$t.Cell(5, 4).Range.Text = "71÷7=10, 1"
$t.Cell(5, 5).Range.Text = "38÷8=4, 6"

$t.Cell(9, 1).Range.Text = "92÷8=11, 4"
$t.Cell(9, 2).Range.Text = "28÷9=3, 1"
$t.Cell(9, 3).Range.Text = "21÷9=2, 3"
$t.Cell(9, 4).Range.Text = "15÷7=2, 1"
$t.Cell(9, 5).Range.Text = "77÷3=25, 2"

$t.Cell(13, 1).Range.Text = "81÷7=11, 4"
$t.Cell(13, 2).Range.Text = "69÷4=17, 1"
$t.Cell(13, 3).Range.Text = "64÷2=32, 0"
$t.Cell(13, 4).Range.Text = "55÷2=27, 1"
$t.Cell(13, 5).Range.Text = "16÷5=3, 1"

$t.Cell(17, 1).Range.Text = "88÷9=9, 7"
$t.Cell(17, 2).Range.Text = "15÷6=2, 3"
$t.Cell(17, 3).Range.Text = "29÷8=3, 5"
$t.Cell(17, 4).Range.Text = "20÷3=6, 2"
$t.Cell(17, 5).Range.Text = "98÷4=24, 2"
